$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    # Force the literal string into the cell without Excel re-interpreting
    # numeric-looking text (e.g. "0.998") as a number: apply a Text number
    # format just long enough to assign the value, then clear formatting so
    # the cell is left with no explicit style (matches the source workbook).
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

$ws.Range("D2").Value = '68.172.83'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '3.828.22'
$ws.Range("E3").Value = '  -1.56%  '
Set-TextValue $ws.Range("D4") '0.998'
$ws.Range("E4").Value = '  -0.12%  '
Set-TextValue $ws.Range("D5") '600.73'
$ws.Range("E5").Value = '  +0.08%  '
Set-TextValue $ws.Range("D6") '169.53'
$ws.Range("E6").Value = '  -0.99%  '
$ws.Range("D7").Value = '3.834.60'
$ws.Range("E7").Value = '  -1.41%  '
$ws.Range("E8").Value = '  +0.02%  '
Set-TextValue $ws.Range("D9") '0.531'
$ws.Range("E9").Value = '  +0.01%  '
Set-TextValue $ws.Range("D10") '0.166'
$ws.Range("E10").Value = '  +0.73%  '
Set-TextValue $ws.Range("D11") '6.52'
$ws.Range("E11").Value = '  +1.79%  '
Set-TextValue $ws.Range("D12") '0.463'
$ws.Range("E12").Value = '  +0.96%  '
Set-TextValue $ws.Range("D13") '0.0000275'
$ws.Range("E13").Value = '  +5.83%  '
Set-TextValue $ws.Range("D14") '37.15'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").Value = '4.471.14'
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("D16").Value = '3.829.08'
$ws.Range("E16").Value = '  -1.78%  '
Set-TextValue $ws.Range("D17") '19.06'
$ws.Range("E17").Value = '  +5.17%  '
$ws.Range("D18").Value = '68.132.37'
$ws.Range("E18").Value = '  -0.31%  '
Set-TextValue $ws.Range("D19") '7.38'
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("E20").Value = '  +0.71%  '
Set-TextValue $ws.Range("D21") '10.75'
$ws.Range("E21").Value = '  -1.07%  '
Set-TextValue $ws.Range("D22") '470.47'
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("E23").Value = '  +0.61%  '
Set-TextValue $ws.Range("D24") '0.0000152'
$ws.Range("E24").Value = '  -4.89%  '
Set-TextValue $ws.Range("D25") '83.88'
$ws.Range("E25").Value = '  +0.51%  '
Set-TextValue $ws.Range("D26") '2.29'
$ws.Range("E26").Value = '  +2.66%  '
Set-TextValue $ws.Range("D27") '12.26'
$ws.Range("E27").Value = '  +1.28%  '
Set-TextValue $ws.Range("D28") '10.38'
$ws.Range("E28").Value = '  +4.32%  '
$ws.Range("E29").Value = '  +0.01%  '
Set-TextValue $ws.Range("D30") '2.94'
$ws.Range("E30").Value = '  -0.45%  '
$ws.Range("D31").Value = '3.971.08'
$ws.Range("E31").Value = '  -1.70%  '
Set-TextValue $ws.Range("D32") '7.73'
$ws.Range("E32").Value = '  -0.25%  '
Set-TextValue $ws.Range("D33") '2.29'
$ws.Range("E33").Value = '  -0.90%  '
Set-TextValue $ws.Range("D34") '30.81'
$ws.Range("E34").Value = '  -1.54%  '
Set-TextValue $ws.Range("D35") '9.31'
$ws.Range("E35").Value = '  -0.95%  '
$ws.Range("D36").Value = '3.791.98'
$ws.Range("E36").Value = '  -1.72%  '
Set-TextValue $ws.Range("D37") '3.81'
$ws.Range("E37").Value = '  +3.38%  '
$ws.Range("E38").Value = '  +1.38%  '
Set-TextValue $ws.Range("D39") '5.98'
$ws.Range("E39").Value = '  +1.31%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D40") '1.02'
$ws.Range("E40").Value = '  -0.75%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D41") '0.139'
$ws.Range("E41").Value = '  -0.89%  '
Set-TextValue $ws.Range("D42") '0.997'
$ws.Range("E42").Value = '  -0.21%  '
Set-TextValue $ws.Range("D43") '0.320'
$ws.Range("E43").Value = '  +2.42%  '
$ws.Range("E44").Value = '  -0.01%  '
Set-TextValue $ws.Range("D45") '1.98'
$ws.Range("E45").Value = '  -0.06%  '
Set-TextValue $ws.Range("D46") '8.78'
$ws.Range("E46").Value = '  +2.00%  '
Set-TextValue $ws.Range("D47") '410.99'
$ws.Range("E47").Value = '  -2.85%  '
Set-TextValue $ws.Range("D48") '46.58'
$ws.Range("E48").Value = '  -1.36%  '
Set-TextValue $ws.Range("D49") '0.000283'
$ws.Range("E49").Value = '  -6.53%  '
Set-TextValue $ws.Range("D50") '142.34'
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("E51").Value = '  +0.26%  '
